$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $old"
    }
}

# Paragraph: "Programa" (Portuguese) - insert line breaks between the three sentences
Replace-Text "integração termo a termo. Equações diferenciais ordinárias" "integração termo a termo. ^lEquações diferenciais ordinárias"
Replace-Text "de 1ª e 2ª ordem.•Séries de Fourier" "de 1ª e 2ª ordem.^l•^lSéries de Fourier"

# Paragraph: "Programa" (English, italic) - insert line breaks between the three sentences
Replace-Text "integration term to term.•First and second order" "integration term to term.^l•First and second order"
Replace-Text "second order differential equations.•Fourier series" "second order differential equations.^l•Fourier series"

# Paragraph: "Bibliografia" - insert line breaks between the five numbered references
Replace-Text "Científicos, 1987.2.BRANNAN" "Científicos, 1987.^l2.BRANNAN"
Replace-Text "LTC ED., 2008.3.ZILL" "LTC ED., 2008.^l3.ZILL"
Replace-Text "Pearson Makron Books2006., v.1 e 2.4.W. Kaplan" "Pearson Makron Books2006., v.1 e 2.^l4.W. Kaplan"
Replace-Text "São Paulo, 1972.5.BOYCE" "São Paulo, 1972.^l5.BOYCE"
